# Bugfix in regional definition: RSAF and SAF were swapped in the
# "Country to region" capital definition (column C = IMAGE region code,
# column D = IMAGE region name). Countries that were tagged 10/"South Africa"
# should be 26/"Rest of South Africa", and South Africa itself (ZAF, row 40),
# which was wrongly tagged 26/"Rest of South Africa", should be 10/"South Africa".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country to region")

# Rows that were 10 / South Africa -> should become 26 / Rest of South Africa
$rowsToRestOfSouthAfrica = @(2, 4, 24, 27, 31, 32, 42, 43, 46, 47, 155)
foreach ($r in $rowsToRestOfSouthAfrica) {
    $ws.Range("C$r").Value = 26
    $ws.Range("D$r").Value = "Rest of South Africa"
}

# Row 40 (South Africa / ZAF) was 26 / Rest of South Africa -> should become 10 / South Africa
$ws.Range("C40").Value = 10
$ws.Range("D40").Value = "South Africa"

# Leftover scratch reference values the author typed while checking the
# correct codes (columns H/I, rows 258-259).
$ws.Range("H258").Value = 10
$ws.Range("I258").Value = "South Africa"
$ws.Range("H259").Value = 26
$ws.Range("I259").Value = "Rest of South Africa"

# Turn on AutoFilter over the table range, and register the corresponding
# hidden workbook-scoped _FilterDatabase defined name (as Excel does).
$ws.Range("A1:D240").AutoFilter() | Out-Null
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "='Country to region'!`$A`$1:`$D`$240")
$fdb.Visible = $false

# Restore the selection that was active when the workbook was saved.
$ws.Range("G244").Select() | Out-Null
